# Disaggregation of commodity Copper
# 1. Rename the commodity "Copper ores and concentrates" to "Copper" on every
#    sheet (each year-sheet has the commodity name in C4).
# 2. A handful of sheets had their D4 cached totals nudged by a single ULP as
#    a side effect of the recalculation that happened when the data was
#    disaggregated.

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)
    if ($ws.Range("C4").Value2 -eq "Copper ores and concentrates") {
        $ws.Range("C4").Value = "Copper"
    }
}

# Tiny last-digit precision updates to the cached D4 values on the affected
# year sheets (sheet index == year - 1999, per xl/_rels/workbook.xml.rels).
$d4Updates = @{
    22 = 55258.16468093192
    24 = 90949.05693803652
    33 = 222769.1412828042
    42 = 827792.0901816025
    45 = 1649624.326648425
    75 = 3220268.906713158
    86 = 2751427.116686261
    91 = 3358437.817936322
    92 = 3395075.375635045
}

foreach ($idx in $d4Updates.Keys) {
    $ws = $wb.Worksheets.Item([int]$idx)
    $ws.Range("D4").Value = $d4Updates[$idx]
}
